$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the parts-list rows with the final, sorted data ---------------
# (Rows are written already in their final alphabetical-by-Device order so
# the subsequent Sort call below is a structural no-op that only records
# the sortState Excel leaves behind after a manual Data > Sort.)

$rows = @(
    @{ Id=1;  Device="GearBox";                  Part=$null;                         Price=$null;   Qty=1; Source="Zepler Stores" },
    @{ Id=2;  Device="I2C Mux";                   Part="PCA9542A";                    Price=0.81;    Qty=1; Source="Farnell" },
    @{ Id=3;  Device="Micro SD Card";              Part=$null;                         Price=4;       Qty=1; Source="Amazon" },
    @{ Id=4;  Device="Micro SD Card Connector";    Part=$null;                         Price=2.04;    Qty=1; Source="Farnell" },
    @{ Id=5;  Device="Microcontroller";             Part="AT32UC3C0512C";               Price=15.39;   Qty=1; Source="Farnell" },
    @{ Id=6;  Device="Motor Driver";                Part=$null;                         Price=1.07;    Qty=2; Source="Farnell" },
    @{ Id=7;  Device="Motors";                      Part=$null;                         Price=$null;   Qty=2; Source="Zepler Stores" },
    @{ Id=8;  Device="OV7670 + Buffer";             Part=$null;                         Price=$null;   Qty=2; Source="Zepler Stores" },
    @{ Id=9;  Device="PCB";                         Part=$null;                         Price=205.48;  Qty=1; Source="PCB Cart" },
    @{ Id=10; Device="Robot Base";                  Part=$null;                         Price=$null;   Qty=1; Source="Zepler Stores" },
    @{ Id=11; Device="SDRAM";                       Part="MICRON - MT48LC4M16A2P-75";    Price=3.24;    Qty=1; Source="Farnell" },
    @{ Id=12; Device="Voltage Regulator";           Part=$null;                         Price=1.03;    Qty=1; Source="Farnell" }
)

$r = 2
foreach ($row in $rows) {
    # Every column is written explicitly (clearing with $null where the new
    # row has no value) so nothing from the pre-edit layout leaks through.
    $ws.Cells.Item($r, 1).Value = $row.Id
    $ws.Cells.Item($r, 2).Value = $row.Device
    $ws.Cells.Item($r, 3).Value = $row.Part
    $ws.Cells.Item($r, 5).Value = $row.Qty
    $ws.Cells.Item($r, 7).Value = $row.Source

    if ($null -ne $row.Price) {
        $ws.Cells.Item($r, 4).Value = $row.Price
        $ws.Cells.Item($r, 6).Formula = "=E$r*D$r"
        $ws.Cells.Item($r, 6).NumberFormat = """£""#,##0.00"
    } else {
        $ws.Cells.Item($r, 4).Value = $null
        $ws.Cells.Item($r, 6).Value = $null
    }
    $r++
}

# --- Record a Data > Sort (by Device, ascending) over B1:G13 ---------------
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("B1")) | Out-Null
$sort.SetRange($ws.Range("B1:G13"))
$sort.Header = 2
$sort.Apply()

# --- Totals row --------------------------------------------------------
$ws.Range("E15").Value = "Total:"
$ws.Range("F15").Formula = "=SUM(F2:F13)"
$ws.Range("F15").NumberFormat = """£""#,##0.00"

# --- Column widths (match the new bestFit sizing as closely as the COM
#     ColumnWidth setter's quantisation allows) ---
$ws.Columns.Item(2).ColumnWidth = 22.3
$ws.Columns.Item(4).ColumnWidth = 19.7
$ws.Columns.Item(7).ColumnWidth = 11.9

# --- Page setup ----------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection -------------------------------------------------------------
$ws.Range("G15").Select()
